{"js": "const replacements = [\n  [\"2026-01-27 Tuesday\", \"2026-01-28 Wednesday\"],\n  [\"442\u00d72=\", \"154\u00d76=\"],\n  [\"146\u00d72=\", \"187\u00d74=\"],\n  [\"397\u00d73=\", \"259\u00d78=\"],\n  [\"867\u00d76=\", \"403\u00d78=\"],\n  [\"303\u00d73=\", \"421\u00d78=\"],\n  [\"860\u00d75=\", \"527\u00d72=\"],\n  [\"839\u00d79=\", \"725\u00d77=\"],\n  [\"573\u00d76=\", \"184\u00d77=\"],\n  [\"653\u00d78=\", \"281\u00d73=\"],\n  [\"369\u00d78=\", \"831\u00d73=\"],\n  [\"819\u00d78=\", \"533\u00d75=\"],\n  [\"274\u00d73=\", \"672\u00d74=\"],\n  [\"329\u00d73=\", \"106\u00d79=\"],\n  [\"822\u00d79=\", \"244\u00d74=\"],\n  [\"104\u00d77=\", \"798\u00d73=\"],\n  [\"374\u00d78=\", \"675\u00d76=\"],\n  [\"345\u00d79=\", \"581\u00d75=\"],\n  [\"127\u00d73=\", \"845\u00d76=\"],\n  [\"750\u00d72=\", \"253\u00d73=\"],\n  [\"637\u00d72=\", \"804\u00d77=\"],\n  [\"609\u00d75=\", \"381\u00d75=\"],\n  [\"991\u00d76=\", \"314\u00d79=\"],\n  [\"700\u00d72=\", \"656\u00d79=\"],\n  [\"203\u00d75=\", \"585\u00d73=\"],\n  [\"874\u00d79=\", \"609\u00d74=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('2026-01-27 Tuesday', '2026-01-28 Wednesday')\n    ,@('442\u00d72=', '154\u00d76=')\n    ,@('146\u00d72=', '187\u00d74=')\n    ,@('397\u00d73=', '259\u00d78=')\n    ,@('867\u00d76=', '403\u00d78=')\n    ,@('303\u00d73=', '421\u00d78=')\n    ,@('860\u00d75=', '527\u00d72=')\n    ,@('839\u00d79=', '725\u00d77=')\n    ,@('573\u00d76=', '184\u00d77=')\n    ,@('653\u00d78=', '281\u00d73=')\n    ,@('369\u00d78=', '831\u00d73=')\n    ,@('819\u00d78=', '533\u00d75=')\n    ,@('274\u00d73=', '672\u00d74=')\n    ,@('329\u00d73=', '106\u00d79=')\n    ,@('822\u00d79=', '244\u00d74=')\n    ,@('104\u00d77=', '798\u00d73=')\n    ,@('374\u00d78=', '675\u00d76=')\n    ,@('345\u00d79=', '581\u00d75=')\n    ,@('127\u00d73=', '845\u00d76=')\n    ,@('750\u00d72=', '253\u00d73=')\n    ,@('637\u00d72=', '804\u00d77=')\n    ,@('609\u00d75=', '381\u00d75=')\n    ,@('991\u00d76=', '314\u00d79=')\n    ,@('700\u00d72=', '656\u00d79=')\n    ,@('203\u00d75=', '585\u00d73=')\n    ,@('874\u00d79=', '609\u00d74=')\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $true, $true, $false, $null, $null, $true, $null, $null, $null, 2)\n}"}
